$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Disponibilizar reserva para subaluguer (" + "utilizador para API
#    central" + ")"  ->  merge the three runs into a single run that
#    carries the identical, already-concatenated text.
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Disponibilizar reserva para subaluguer \(utilizador para API central\)", `
    $true, $false, $true, $false, $false, $true, 1, $false, `
    "Disponibilizar reserva para subaluguer (utilizador para API central)", 2)

# ------------------------------------------------------------------
# 2) "Realizar pedido da disponibilidade de lugar (" + "utilizador
#    para API central" + ")"  ->  merge the three runs into one.
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Realizar pedido da disponibilidade de lugar \(utilizador para API central\)", `
    $true, $false, $true, $false, $false, $true, 1, $false, `
    "Realizar pedido da disponibilidade de lugar (utilizador para API central)", 2)

# ------------------------------------------------------------------
# 3) "Obter valor a pagar de certo lugar (utilizador) " is rewritten
#    as "Realizar pedido para valor a pagar de certo lugar de
#    estacionamento (utilizador) " split across four runs. We rebuild
#    the paragraph via its OOXML so the four <w:r> elements survive
#    the save instead of being re-coalesced into one.
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("Obter valor a pagar de certo lugar \(utilizador\) ", `
    $true, $false, $true, $false, $false, $true, 1, $false, "", 0)

if ($target.Find.Found) {
    $p3 = $target.Paragraphs(1).Range
    $pkg = $p3.XML()

    # Pull the paragraph's own opening <w:p ...> tag (paraId/textId/rsid
    # bookkeeping attributes) straight out of the live XML instead of
    # hard-coding it, so we only hard-code the formatting that the diff
    # actually pins down (style + numbering + spacing).
    $bodyStart = $pkg.IndexOf("<w:body>")
    $paraStart = $pkg.IndexOf("<w:p", $bodyStart)
    $paraOpenEnd = $pkg.IndexOf(">", $paraStart) + 1
    $openTag = $pkg.Substring($paraStart, $paraOpenEnd - $paraStart)
    $paraEnd = $pkg.IndexOf("</w:p>", $paraStart) + 6

    $pPr = '<w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr>'

    $runs = '<w:r><w:t>Realizar pedido para</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve"> valor a pagar de certo lugar</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve"> de estacionamento</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve"> (utilizador) </w:t></w:r>'

    $newParaXml = $openTag + $pPr + $runs + '</w:p>'

    $newPkg = $pkg.Substring(0, $paraStart) + $newParaXml + $pkg.Substring($paraEnd)

    $p3.InsertXML($newPkg)
}
